$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.493.77'
$ws.Range('E2').Value = '  +2.25%  '
$ws.Range('D3').Value = '1.873.02'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.35'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.014'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4783'
$ws.Range('E7').Value = '  +0.76%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3774'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07376'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9377'
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.75'
$ws.Range('E11').Value = '  +5.83%  '
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').Value = '1.874.98'
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.590'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.92'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.015'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008920'
$ws.Range('E18').Value = '  +3.34%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.014'
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.92'
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').Value = '27.527.18'
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.134'
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.961'
$ws.Range('E24').Value = '  +1.91%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.03'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '115.97'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.997'
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.08935'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.346'
$ws.Range('E31').Value = '  +1.77%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.221'
$ws.Range('E32').Value = '  +4.34%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.612'
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7534'
$ws.Range('E34').Value = '  +0.75%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.692'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02060'
$ws.Range('E36').Value = '  +5.86%  '
$ws.Range('E37').Value = '  +2.67%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05305'
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5359'
$ws.Range('E40').Value = '  +2.93%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.090'
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1528'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.428'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4833'
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.60'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.015'
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.663'
$ws.Range('E47').Value = '  +3.79%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.07'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '67.23'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.8951'
$ws.Range('E51').Value = '  +1.27%  '
